$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header text change: "model" -> "Model"
$ws.Range("A1").Value = "Model"

# Updated metric values (rows 2-6, columns B-E)
$ws.Range("B2").Value = 0.00144795561293999
$ws.Range("C2").Value = 0.03805201194339124
$ws.Range("D2").Value = 0.02670462342758718
$ws.Range("E2").Value = 0.682346721349905

$ws.Range("B3").Value = 0.001464577242493714
$ws.Range("C3").Value = 0.03826979543313126
$ws.Range("D3").Value = 0.02683823725879792
$ws.Range("E3").Value = 0.6857607728492968

$ws.Range("B4").Value = 0.001463096896369838
$ws.Range("C4").Value = 0.0382504496231069
$ws.Range("D4").Value = 0.02710021199128784
$ws.Range("E4").Value = 0.6924546549134173

$ws.Range("B5").Value = 0.001477500091158331
$ws.Range("C5").Value = 0.03843826337334104
$ws.Range("D5").Value = 0.02724523860907797
$ws.Range("E5").Value = 0.6961603217402085

$ws.Range("B6").Value = 0.002915035194714119
$ws.Range("C6").Value = 0.05399106587866292
$ws.Range("D6").Value = 0.03913644279664259
$ws.Range("E6").Value = 1

# Matches the final selection state recorded in the sheet view
[void]$ws.Range("E13").Select()
